$d = $word.ActiveDocument

# 1) "... for None." -> "... for Arraignment."
$d.Content.Find.Execute(" for None", $false, $false, $false, $false, $false, $true, 1, $false, " for Arraignment", 2) | Out-Null

# 2) Heading "None Conditions" -> "Recognizance (OR) Bond Conditions"
$d.Content.Find.Execute("None", $true, $false, $false, $false, $false, $true, 1, $false, "Recognizance (OR) Bond", 2) | Out-Null

# 3) Insert a new list paragraph ("The defendant shall execute a personal
#    recognizance bond.") right before the "Defendant shall behave lawfully..."
#    list item under "Non-Financial Conditions of Release:".
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "Defendant shall behave lawfully*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the 'Defendant shall behave lawfully...' paragraph."
}

$target = $d.Paragraphs($targetIndex)
$target.Range.InsertParagraphBefore() | Out-Null

# The freshly-inserted (empty) paragraph now occupies $targetIndex and
# already inherited the list formatting (ListParagraph style / numPr /
# tabs / jc) from the paragraph that originally sat there.
$newPara = $d.Paragraphs($targetIndex)

$sentence = "The defendant shall execute a personal recognizance bond"
$newPara.Range.InsertAfter($sentence + ".")

$fullRange = $newPara.Range
$underlineRange = $d.Range($fullRange.Start, $fullRange.Start + $sentence.Length)
$underlineRange.Font.Underline = 1
